# edit.ps1 - reproduce the authoring change described by the diff:
#   1. Slide 6's table gets a new table style (tableStyleId GUID change).
#   2. The theme actually driving the slide master's look is switched from
#      the "Integral" palette to the (built-in) "Office Theme" palette.
#
# Helper to build the little-endian 0x00BBGGRR integer that PowerPoint's
# RGB() / ColorFormat.RGB property expects from a "RRGGBB" hex string.
function Hex-ToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Table on slide 6 ("SOURCES OF FINANCE"): apply the new table style
# ------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{2F2769BA-216D-473A-9E1F-179C942E8602}")
    }
}

# ------------------------------------------------------------------
# 2) Re-colour the presentation's active theme (the one bound to the
#    slide master) from "Integral" to the standard "Office Theme"
#    12-colour palette.
# ------------------------------------------------------------------
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Hex-ToRgbInt($officeColors[$i - 1])
}

Write-Host "Applied table style + theme colour updates"
